$p = $ppt.ActivePresentation

# --- Slide 1: title slide text -> new Jinja-style template formatters ----
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Characters().Text = "Client Overview: {{ client_name | format_string('title') }}"
$s1.Shapes.Item(2).TextFrame.TextRange.Characters().Text = "Generated on: {{ contract_date | format_date('medium') }}"

# --- Slide 2: "Risk Assessment" -> "Financial & Status Highlights" -------
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Characters().Text = "Financial & Status Highlights"

$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange
$tr2.Text = "a`rb`rc"
$tr2.Paragraphs(2).IndentLevel = 2
$tr2.Paragraphs(3).IndentLevel = 2
$tr2.Paragraphs(1).Characters().Text = "Key Metrics:"
$tr2.Paragraphs(2).Characters().Text = "Contract Value (BRL): {{ contract_value | format_currency('BRL') }}"
$tr2.Paragraphs(3).Characters().Text = "Current Status: {{ status_code | format_logic('10=Green (Go)', '20=Yellow (Hold)', 'Red (Stop)') }}"

# Remove the old "Confidential - DocGenius Generated" text box entirely.
[void]$s2.Shapes.Item(3).Delete()

# --- Slide 3 (new): "Audit Checkpoints (Boolean Formats)" -----------------
$s3 = $p.Slides.Add(3, 2)
$s3.Shapes.Item(1).TextFrame.TextRange.Characters().Text = "Audit Checkpoints (Boolean Formats)"

$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange
$tr3.Text = "a`rb"
$tr3.Paragraphs(1).Characters().Text = "Is Active User? -> {{ is_active | format_bool('truefalse') }}"
$tr3.Paragraphs(2).Characters().Text = "Debt Clearance Checkbox: [ {{ has_debt | format_bool('checkbox') }} ]"
[void]$tr3.InsertBefore("`r")
